$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Tear down the old "name + verb" table (A1:A4 names, C1:D14 formulas) but
#    keep the verb-stem list (B1:B5) alive so its shared-string ids survive,
#    then relocate it one row down and three columns right (B1:B5 -> E2:E6).
# ---------------------------------------------------------------------------
$ws.Range("A1:A4").Clear()
$ws.Range("C1:D14").Clear()
$ws.Range("B1:B5").Cut($ws.Range("E2:E6"))

# ---------------------------------------------------------------------------
# 2) New shape names for agent/patient (columns B & D, rows 2-4)
# ---------------------------------------------------------------------------
$ws.Range("B2").Value = "Square"
$ws.Range("B3").Value = "Triangle"
$ws.Range("B4").Value = "Star"
$ws.Range("D2").Value = "Triangle"
$ws.Range("D3").Value = "Star"
$ws.Range("D4").Value = "Heart"

# ---------------------------------------------------------------------------
# 3) Header row
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "AgentName"
$ws.Range("B1").Value = "AgentShape"
$ws.Range("C1").Value = "PatientName"
$ws.Range("D1").Value = "PatientShape"

# ---------------------------------------------------------------------------
# 4) Character names (columns A & C, rows 2-4)
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = "Kyle"
$ws.Range("A3").Value = "Lily"
$ws.Range("A4").Value = "Zach"
$ws.Range("C2").Value = "Lily"
$ws.Range("C3").Value = "Zach"
$ws.Range("C4").Value = "Melissa"

# ---------------------------------------------------------------------------
# 5) Remaining headers (G, E, F, H, I, J)
# ---------------------------------------------------------------------------
$ws.Range("G1").Value = "ProgressiveSentence"
$ws.Range("E1").Value = "VerbStem"
$ws.Range("F1").Value = "Verb"
$ws.Range("H1").Value = "PastSentence"
$ws.Range("I1").Value = "ShapeProgressiveSentence"
$ws.Range("J1").Value = "ShapePastSentence"

# ---------------------------------------------------------------------------
# 6) Sentence-builder formulas
# ---------------------------------------------------------------------------
$ws.Range("G2").Formula = '=CONCATENATE($A2, " ", $B2, " is ",$E2, "ing ",$C2, " ", $D2,".")'
$ws.Range("I2").Formula = '=CONCATENATE("The ", $B2, " is ",$E2, "ing ","the ", $D2,".")'

$ws.Range("G3:G8").Formula = '=CONCATENATE($A3, " ", $B3, " is ",E3, "ing ",C3, " ", D3,".")'

$ws.Range("G9").Formula = '=CONCATENATE(B4," is ",E2, "ing ",D4,".")'
$ws.Range("G10").Formula = '=CONCATENATE(D4," is ",E2, "ing ",B2,".")'
$ws.Range("G11").Formula = '=CONCATENATE(D4," is ",E2, "ing ",B3,".")'
$ws.Range("G12").Formula = '=CONCATENATE(D4," is ",E2, "ing ",B4,".")'

# G15 keeps the wrap-style formatting from the G3:G8 fill but no content
$ws.Range("G15").Value = ""

# ---------------------------------------------------------------------------
# 7) Unrelated scratch counter table (A13:B17)
# ---------------------------------------------------------------------------
$ws.Range("A13").Value = 1
$ws.Range("A14").Value = 2
$ws.Range("A15").Value = 3
$ws.Range("A16").Value = 4
$ws.Range("A17").Value = 5

$ws.Range("B13").Formula = "=A13+2"
$ws.Range("B14:B17").Formula = "=A14+2"

# ---------------------------------------------------------------------------
# 8) Styling: the whole ProgressiveSentence fill column is left-aligned (this
#    reuses the pre-existing "left align" style), then the Kyle/Square demo
#    row gets a yellow highlight on top (forces a new merged fill+align style
#    for G2/I2, and a fill-only style for A2:E2).
# ---------------------------------------------------------------------------
$ws.Range("G2").HorizontalAlignment = -4131
$ws.Range("I2").HorizontalAlignment = -4131
$ws.Range("G3:G8").HorizontalAlignment = -4131
$ws.Range("G15").HorizontalAlignment = -4131

$ws.Range("A2:E2").Interior.Color = 65535
$ws.Range("G2").Interior.Color = 65535
$ws.Range("I2").Interior.Color = 65535

# ---------------------------------------------------------------------------
# 9) Column widths / layout
# ---------------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 12.5
$ws.Columns.Item(3).ColumnWidth = 12.5
$ws.Columns.Item(4).ColumnWidth = 12.5
$ws.Columns.Item(7).ColumnWidth = 36.5
$ws.Columns.Item(8).ColumnWidth = 34

# ---------------------------------------------------------------------------
# 10) Sheet view: scroll / selection
# ---------------------------------------------------------------------------
$ws.Range("E12").Select()
